$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 195 and 196, shifting the existing rows
# (old 195..220) down to (197..222).
$ws.Range("A195:A196").EntireRow.Insert()

# Populate new row 195
$ws.Cells.Item(195, 1).Value = 3
$ws.Cells.Item(195, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(195, 3).Value = "Coquimbo"
$ws.Cells.Item(195, 4).Value = 44491
$ws.Cells.Item(195, 5).Value = 5
$ws.Cells.Item(195, 6).Value = 100112013
$ws.Cells.Item(195, 7).Value = "Alcachofa"
$ws.Cells.Item(195, 8).Value = "Española"
$ws.Cells.Item(195, 9).Value = "Extra"
$ws.Cells.Item(195, 10).Value = 13300
$ws.Cells.Item(195, 11).Value = 370
$ws.Cells.Item(195, 12).Value = 380
$ws.Cells.Item(195, 13).Value = 375
$ws.Cells.Item(195, 14).Value = "`$/unidad"
$ws.Cells.Item(195, 15).Value = "Llay Llay"
$ws.Cells.Item(195, 16).Value = 375
$ws.Cells.Item(195, 17).Value = 1
$ws.Cells.Item(195, 18).Value = "Hortaliza"

# Populate new row 196
$ws.Cells.Item(196, 1).Value = 3
$ws.Cells.Item(196, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(196, 3).Value = "Coquimbo"
$ws.Cells.Item(196, 4).Value = 44491
$ws.Cells.Item(196, 5).Value = 5
$ws.Cells.Item(196, 6).Value = 100112013
$ws.Cells.Item(196, 7).Value = "Alcachofa"
$ws.Cells.Item(196, 8).Value = "Española"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 6500
$ws.Cells.Item(196, 11).Value = 300
$ws.Cells.Item(196, 12).Value = 300
$ws.Cells.Item(196, 13).Value = 300
$ws.Cells.Item(196, 14).Value = "`$/unidad"
$ws.Cells.Item(196, 15).Value = "Llay Llay"
$ws.Cells.Item(196, 16).Value = 300
$ws.Cells.Item(196, 17).Value = 1
$ws.Cells.Item(196, 18).Value = "Hortaliza"
